# ---------------------------------------------------------------------------
# Update with Correct Forecast output
#
# Sheet1 ("Sheet1" -> "Sales vs PO"):
#   - insert a new column C ("Order Week") holding the previous week's order
#     date (i.e. the old "ds" values), shifting the old PO_Requested_Qty
#     column to D
#   - bump every "ds" date in column A forward by one week (+6, since the
#     values are already 1 day into the new week)
#   - the new PO_Requested_Qty column (D) is reset to 0 for every row
#     (forecast placeholder column)
# Add three new sheets with the forecast/analysis output:
#   - "Weekly Growth"     : ds / PO_Requested_Qty / Growth%
#   - "Volume Insights"   : Total / Average / Max / Min PO quantity
#   - "Prediction Info"   : Predicted_Next_Week_PO_Quantity
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Sales vs PO"
# ------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item(1)
$wsSales.Name = "Sales vs PO"

# Insert a fresh column C; this shifts the existing C column (header
# "PO_Requested_Qty" + its values) into column D, formatting included.
$wsSales.Columns("C:C").Insert()

# Header for the newly inserted column.
$wsSales.Cells.Item(1, 3).Value = "Order Week"

# Give the new column's data cells the same (date) number format style as
# column A by copying formats from A2:A26 onto C2:C26.
$wsSales.Range("A2:A26").Copy()
$wsSales.Range("C2:C26").PasteSpecial(-4122)

# Original "ds" values (column A) before this edit - row 2 .. row 26.
$origDs = @(45481,45488,45495,45502,45509,45516,45523,45530,45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)

for ($i = 0; $i -lt $origDs.Length; $i++) {
    $r = $i + 2
    $wsSales.Cells.Item($r, 1).Value = $origDs[$i] + 6   # ds moves one week later
    $wsSales.Cells.Item($r, 3).Value = $origDs[$i]        # Order Week = old ds
    $wsSales.Cells.Item($r, 4).Value = 0                  # PO_Requested_Qty placeholder
}

# ------------------------------------------------------------------
# Sheet 2: "Weekly Growth"
# ------------------------------------------------------------------
$wsGrowth = $wb.Worksheets.Add($null, $wsSales)
$wsGrowth.Name = "Weekly Growth"

# Copy header + date-column formatting from "Sales vs PO" so the new sheet
# matches its look (bold/bordered header, date-formatted first column).
$wsSales.Range("A1:C1").Copy()
$wsGrowth.Range("A1:C1").PasteSpecial(-4122)
$wsSales.Range("A2:A12").Copy()
$wsGrowth.Range("A2:A12").PasteSpecial(-4122)

$wsGrowth.Cells.Item(1, 1).Value = "ds"
$wsGrowth.Cells.Item(1, 2).Value = "PO_Requested_Qty"
$wsGrowth.Cells.Item(1, 3).Value = "Growth%"

$growthRows = @(
    @(45488, 32, 0),
    @(45495, 16, -50),
    @(45509, 16, 0),
    @(45516, 16, 0),
    @(45523, 32, 100),
    @(45530, 16, -50),
    @(45537, 96, 500),
    @(45544, 32, -66.66666666666667),
    @(45551, 208, 550),
    @(45586, 208, 0),
    @(45607, 16, -92.30769230769231)
)

for ($i = 0; $i -lt $growthRows.Length; $i++) {
    $r = $i + 2
    $row = $growthRows[$i]
    $wsGrowth.Cells.Item($r, 1).Value = $row[0]
    $wsGrowth.Cells.Item($r, 2).Value = $row[1]
    $wsGrowth.Cells.Item($r, 3).Value = $row[2]
}

# ------------------------------------------------------------------
# Sheet 3: "Volume Insights"
# ------------------------------------------------------------------
$wsVolume = $wb.Worksheets.Add($null, $wsGrowth)
$wsVolume.Name = "Volume Insights"

$wsSales.Range("A1:C1").Copy()
$wsVolume.Range("A1:D1").PasteSpecial(-4122)

$wsVolume.Cells.Item(1, 1).Value = "Total_PO_Quantity"
$wsVolume.Cells.Item(1, 2).Value = "Average_PO_Quantity"
$wsVolume.Cells.Item(1, 3).Value = "Max_PO_Quantity"
$wsVolume.Cells.Item(1, 4).Value = "Min_PO_Quantity"

$wsVolume.Cells.Item(2, 1).Value = 688
$wsVolume.Cells.Item(2, 2).Value = 62.54545454545455
$wsVolume.Cells.Item(2, 3).Value = 208
$wsVolume.Cells.Item(2, 4).Value = 16

# ------------------------------------------------------------------
# Sheet 4: "Prediction Info"
# ------------------------------------------------------------------
$wsPredict = $wb.Worksheets.Add($null, $wsVolume)
$wsPredict.Name = "Prediction Info"

$wsSales.Range("A1").Copy()
$wsPredict.Range("A1").PasteSpecial(-4122)

$wsPredict.Cells.Item(1, 1).Value = "Predicted_Next_Week_PO_Quantity"
$wsPredict.Cells.Item(2, 1).Value = 136.7272727272727

# ------------------------------------------------------------------
# Restore the first sheet as the active / selected tab.
# ------------------------------------------------------------------
$wsSales.Activate()
